$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '307.94'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-2.05%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.85'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.67%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.039'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.86%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07636'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-3.21%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.239'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-2.80%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.607'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-3.78%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.452'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-4.66%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9088'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.56%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1019'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-8.10%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1774'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.99%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09100'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.05%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04396'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.51%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1052'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.31%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001271'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.81%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005810'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.06%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.367'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.39%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.3299'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.39%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-6.92%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1358'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.21%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2719'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.27%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04154'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.59%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-3.52%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004092'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.87%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001298'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '5.17%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003002'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.61%'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-1.84%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05170'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-3.38%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007765'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.97%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1309'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-4.03%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007077'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-6.70%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001945'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-6.36%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008029'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-3.20%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3063'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.94%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006363'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-6.75%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-1.53%'
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.004393'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '5.66%'
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004871'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '42.77%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002097'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-1.53%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001997'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-1.53%'
